$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.830.91"
$ws.Range("E2").Value = "  -0.35%  "

$ws.Range("D3").Value = "1.895.72"
$ws.Range("E3").Value = "  +0.13%  "

$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7984"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.55%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3169"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.60%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.54"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07045"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08058"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7718"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.52%  "

$ws.Range("D13").Value = "1.890.44"
$ws.Range("E13").Value = "  -0.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.329"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.28%  "

$ws.Range("D16").Value = "29.813.54"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.982"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007713"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.321"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +20.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.16%  "

$ws.Range("D23").Value = "2.136.66"
$ws.Range("E23").Value = "  -0.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1661"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.337"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.81%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.84%  "

$ws.Range("E28").Value = "  -0.65%  "

$ws.Range("E29").Value = "  -0.59%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.398"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.538"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.434"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05700"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.049"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.261"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.76%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7390"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9985"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.633"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01912"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.17%  "

$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4415"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.818"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8459"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9993"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.18%  "

$ws.Range("D46").Value = "1.036.19"
$ws.Range("E46").Value = "  +4.90%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.876"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.57%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.957"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.445"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.68%  "

$ws.Range("D51").Value = "2.030.49"
$ws.Range("E51").Value = "  -0.72%  "
